$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet gained a new daily data row (23), duplicating the values that
# were previously the last row (22). The "date-only" number format that
# used to sit on C22 moves down onto the new C23, while C22 switches to
# the regular datetime format shared by the rest of the data rows.

# 1. Duplicate row 22 (values + formatting) down into row 23 so C23
#    inherits the date-only style (s=3) that C22 currently has.
$ws.Range("A22:C22").Copy($ws.Range("A23:C23"))

# 2. Re-apply the common datetime style (taken from C2) to C22 and make
#    sure its value is unchanged.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("C22").Value = 45754
